$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new D (Price) value; rows not listed keep their existing D value
$newD = @{
    2  = "22.477.79"
    3  = "1.572.21"
    5  = "1.002"
    6  = "291.82"
    7  = "0.3721"
    8  = "49.86"
    9  = "0.3400"
    10 = "1.147"
    11 = "0.07547"
    12 = "1.002"
    13 = "21.23"
    14 = "6.039"
    15 = "6.959"
    16 = "1.572.36"
    17 = "0.00001123"
    18 = "90.80"
    19 = "0.06760"
    21 = "6.293"
    22 = "16.34"
    24 = "22.470.26"
    25 = "2.377"
    26 = "2.624"
    28 = "149.12"
    29 = "5.046"
    30 = "125.36"
    31 = "1.746.59"
    32 = "1.081"
    33 = "6.206"
    35 = "9.803"
    36 = "0.08353"
    37 = "0.02480"
    38 = "0.2302"
    39 = "1.341"
    40 = "0.06549"
    41 = "5.442"
    43 = "0.6230"
    45 = "14.01"
    46 = "3.814"
    47 = "0.5845"
    48 = "130.73"
    49 = "2.069"
    51 = "0.07327"
}

# Mapping of row -> new E (Volume(1h)) value (inner percentage text, padded with two spaces on each side)
$newE = @{
    2  = "+0.45%"
    3  = "+0.12%"
    4  = "-0.10%"
    5  = "-0.02%"
    6  = "+0.21%"
    7  = "-1.12%"
    8  = "+0.17%"
    9  = "-0.50%"
    10 = "+0.48%"
    11 = "-0.86%"
    12 = "-0.12%"
    13 = "+0.41%"
    14 = "+0.79%"
    15 = "+0.44%"
    16 = "+0.10%"
    17 = "-0.88%"
    18 = "+0.60%"
    19 = "+0.18%"
    20 = "-0.03%"
    21 = "+1.68%"
    22 = "-2.26%"
    23 = "+1.54%"
    24 = "+0.43%"
    25 = "-0.71%"
    26 = "-1.00%"
    27 = "-0.37%"
    28 = "+1.37%"
    29 = "+0.06%"
    30 = "-0.99%"
    31 = "+0.01%"
    32 = "+9.30%"
    33 = "+2.04%"
    34 = "+0.17%"
    35 = "-3.20%"
    36 = "-1.89%"
    37 = "-1.97%"
    38 = "-0.18%"
    39 = "-2.32%"
    40 = "+0.86%"
    41 = "+0.79%"
    42 = "+0.37%"
    43 = "-1.60%"
    44 = "-0.06%"
    45 = "-0.51%"
    46 = "+0.61%"
    47 = "-1.71%"
    48 = "+4.77%"
    49 = "-0.88%"
    50 = "-4.28%"
    51 = "+0.07%"
}

foreach ($row in $newD.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $newD[$row]
}

foreach ($row in $newE.Keys) {
    $cell = $ws.Range("E$row")
    $cell.NumberFormat = "@"
    $cell.Value = "  " + $newE[$row] + "  "
}
